$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove all existing hyperlinks (and their relationships) so we can rebuild them cleanly
$ws.Cells.Hyperlinks.Delete()

# Row 2
$ws.Range("A2").Value = "2026-01-13 18:32:16"
$ws.Range("B2").Value = "大手SIer等のAIソリューション開発・導入を支援してくださるエンジニア・PM募集"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5455098"
$ws.Range("G2").Value = 375
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

# Row 3
$ws.Range("A3").Value = "2026-01-13 18:32:16"
$ws.Range("B3").Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5445159"
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# Row 4
$ws.Range("A4").Value = "2026-01-13 18:32:16"
$ws.Range("B4").Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5445154"
$ws.Range("G4").Value = 368
$ws.Range("H4").Value = "🔥AI,Ai ◆開発"

# Row 5
$ws.Range("A5").Value = "2026-01-13 18:32:16"
$ws.Range("B5").Value = "【募集】Python / Docker 日次データ スクレイピングシステム構築"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5469627"
$ws.Range("G5").Value = 248
$ws.Range("H5").Value = "🔥Python ◆スクレイピング"

# Row 6
$ws.Range("A6").Value = "2026-01-13 18:32:16"
$ws.Range("B6").Value = "【急募】FXツール開発のプロフェッショナルを探しています!"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5470011"
$ws.Range("G6").Value = 123
$ws.Range("H6").Value = "◆ツール,開発"

# Row 7
$ws.Range("A7").Value = "2026-01-13 18:32:16"
$ws.Range("B7").Value = "自社システムの開発・保守エンジニア募集★☆カメラ面談あり"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5469878"
$ws.Range("G7").Value = 75
$ws.Range("H7").Value = "◆開発"

# Row 8
$ws.Range("A8").Value = "2026-01-13 18:32:16"
$ws.Range("B8").Value = "GoogleCloudを利用したアジャイル開発共通基盤のSREエンジニアの募集"
$ws.Range("C8").Value = "システム開発"
$ws.Range("D8").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E8").Value = "期限情報なし"
$ws.Range("F8").Value = "https://www.lancers.jp/work/detail/5457458"
$ws.Range("G8").Value = 75
$ws.Range("H8").Value = "◆開発"

# Row 9
$ws.Range("A9").Value = "2026-01-13 18:32:16"
$ws.Range("B9").Value = "【要日本語ネイティブ】設計書なしレガシーPHP5.3システムの完全再現移行(Win→Linux)"
$ws.Range("C9").Value = "システム開発"
$ws.Range("D9").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E9").Value = "期限情報なし"
$ws.Range("F9").Value = "https://www.lancers.jp/work/detail/5470085"
$ws.Range("G9").Value = 55
$ws.Range("H9").Value = "○PHP"

# Row 10
$ws.Range("A10").Value = "2026-01-13 18:32:16"
$ws.Range("B10").Value = "進行管理およびチームディレクションを担当"
$ws.Range("C10").Value = "システム開発"
$ws.Range("D10").Value = "~ 5,000 円 / 固定"
$ws.Range("E10").Value = "期限情報なし"
$ws.Range("F10").Value = "https://www.lancers.jp/work/detail/5418064"
$ws.Range("G10").Value = 30
$ws.Range("H10").Value = "◇管理"

# Row 11
$ws.Range("A11").Value = "2026-01-13 18:32:16"
$ws.Range("B11").Value = "【緊急】AWS上の稼働中Webサイトを最新GitHubに再構築"
$ws.Range("C11").Value = "システム開発"
$ws.Range("D11").Value = "5,000 円 ~ 10,000 円 / 固定"
$ws.Range("E11").Value = "期限情報なし"
$ws.Range("F11").Value = "https://www.lancers.jp/work/detail/5469840"
$ws.Range("G11").Value = 30
$ws.Range("H11").Value = "◇サイト"

# Row 12
$ws.Range("A12").Value = "2026-01-13 18:32:16"
$ws.Range("B12").Value = "金融機関の入出金伝票印刷システム構築依頼"
$ws.Range("C12").Value = "システム開発"
$ws.Range("D12").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E12").Value = "期限情報なし"
$ws.Range("F12").Value = "https://www.lancers.jp/work/detail/5470403"
$ws.Range("G12").Value = 28

# Row 13
$ws.Range("A13").Value = "2026-01-13 18:32:16"
$ws.Range("B13").Value = "初回 【訪問看護向けオンライン請求システム】エラー解消・仕様確認サポート担当募集"
$ws.Range("C13").Value = "システム開発"
$ws.Range("D13").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E13").Value = "期限情報なし"
$ws.Range("F13").Value = "https://www.lancers.jp/work/detail/5470035"
$ws.Range("G13").Value = 28

# Row 14
$ws.Range("A14").Value = "2026-01-13 18:32:16"
$ws.Range("B14").Value = "コールセンター通話・SMSの発信判断を統合する顧客DB(MUCS)PoC"
$ws.Range("C14").Value = "システム開発"
$ws.Range("D14").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E14").Value = "期限情報なし"
$ws.Range("F14").Value = "https://www.lancers.jp/work/detail/5470137"
$ws.Range("G14").Value = 25

# Row 15
$ws.Range("A15").Value = "2026-01-13 18:32:16"
$ws.Range("B15").Value = "《長期レギュラー》公的機関Web運用の要となる、ディレクター募集"
$ws.Range("C15").Value = "システム開発"
$ws.Range("D15").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E15").Value = "期限情報なし"
$ws.Range("F15").Value = "https://www.lancers.jp/work/detail/5470150"
$ws.Range("G15").Value = 18

# Row 16
$ws.Range("A16").Value = "2026-01-13 18:32:16"
$ws.Range("B16").Value = "フロント実装済み!音楽権利マーケットプレイス「HITOON」のバックエンド・決済機能実装"
$ws.Range("C16").Value = "システム開発"
$ws.Range("D16").Value = "100,000 円 ~ 200,000 円 / 固定"
$ws.Range("E16").Value = "期限情報なし"
$ws.Range("F16").Value = "https://www.lancers.jp/work/detail/5470263"
$ws.Range("G16").Value = 18

# Row 17
$ws.Range("A17").Value = "2026-01-13 18:32:16"
$ws.Range("B17").Value = "TradingViewインジケーター作成"
$ws.Range("C17").Value = "システム開発"
$ws.Range("D17").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E17").Value = "期限情報なし"
$ws.Range("F17").Value = "https://www.lancers.jp/work/detail/5470325"
$ws.Range("G17").Value = 13

# Row 18
$ws.Range("A18").Value = "2026-01-13 18:32:16"
$ws.Range("B18").Value = "n8n 初期構築・セットアップ(小規模/検証用途)"
$ws.Range("C18").Value = "システム開発"
$ws.Range("D18").Value = "10,000 円 ~ 20,000 円 / 固定"
$ws.Range("E18").Value = "期限情報なし"
$ws.Range("F18").Value = "https://www.lancers.jp/work/detail/5469826"
$ws.Range("G18").Value = 10

# Re-create hyperlinks for column F across all data rows, with Hyperlink style
for ($r = 2; $r -le 18; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value()
    $ws.Hyperlinks.Add($cell, $url)
    $cell.Style = "Hyperlink"
}
